# Updates the vm_pu.xlsx "res_bus" results sheet for the 380 kV case.
# Columns B,C,D,E,F,I,J,K,L,M,N (column G stays 1, column H has no data)
# for rows 2-25 are replaced with newly computed voltage magnitudes (p.u.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @"
1.02,1.039150021316851,1.059125404531306,1.037458983051559,1.064343318486736,1.046178437044675,1.044243598274249,1.061855761735471,1.040249807895572,1.067059511524012,1.045726543960213
1.02,1.041997152445915,1.06026309633505,1.039948474213141,1.065956170209472,1.046592424316129,1.046726335336052,1.06280680266365,1.042544671551722,1.068485552901673,1.048212806793418
1.02,1.043827916392179,1.06099358186998,1.041548874624282,1.066992160615648,1.046855848693336,1.04832159276424,1.063415941239325,1.044018887305748,1.06940019573547,1.049810329670133
1.02,1.044594887325086,1.061299337961584,1.042219242872169,1.067425895190763,1.04696553568905,1.048989614784935,1.063670545620072,1.044636145022113,1.069782802142126,1.05047930035872
1.02,1.044723510009375,1.061350597722327,1.042331659451509,1.067498616730894,1.046983890929072,1.049101626649734,1.063713208704255,1.044739640351843,1.069846932332361,1.050591471293215
1.02,1.043838175134322,1.06099767263769,1.041557841618425,1.066997963212587,1.04685731847544,1.048330529130309,1.063419349056368,1.044027144884217,1.069405315603793,1.049819278726867
1.02,1.040114654032177,1.059511079928487,1.038302525426937,1.064889988789598,1.046319273790268,1.045085018573256,1.062178475496173,1.041027624318541,1.06754314460047,1.046569159172557
1.02,1.033461639844988,1.056847156001075,1.032483076390574,1.061115570371905,1.045336629414334,1.039276890642237,1.059943194463023,1.035657169669296,1.06419833732506,1.040752783033675
1.02,1.028959719179649,1.055040125703402,1.028543255723792,1.058556958347916,1.04465768424386,1.035340552441626,1.05841904756932,1.032015709648624,1.061923853260197,1.036810854781314
1.02,1.026993423980588,1.054250000754078,1.026822021362458,1.05743853842849,1.044357892755899,1.033619847777511,1.057750741385517,1.030423486049693,1.060927943803288,1.035087706519208
1.02,1.026260411963452,1.053955335478868,1.026180299552383,1.057021484816313,1.04424565212695,1.032978174470532,1.057501224984147,1.029829661439194,1.06055631986995,1.034445121962533
1.02,1.026417766714164,1.054018595892063,1.026318060095652,1.057111018341069,1.044269768368006,1.033115931339588,1.057554805323536,1.029957148850635,1.060636111991754,1.034583074462144
1.02,1.026932887462084,1.054225667805472,1.026769025464167,1.057404097993234,1.044348633033919,1.033566858920398,1.057730142492616,1.030374449796804,1.060897260142671,1.035034642411841
1.02,1.027249917221732,1.054353094813986,1.027046562319276,1.05758445775922,1.04439710653438,1.033844352809326,1.057838003447716,1.030631241948349,1.061057935818725,1.03531253037392
1.02,1.029089850608596,1.055092399957751,1.028657159523217,1.058630958946533,1.044677457057757,1.035454400368494,1.058463223036413,1.032121047955591,1.061989712463483,1.036924864385297
1.02,1.030239389969295,1.055554075576704,1.029663298030007,1.05928455495539,1.044851750809584,1.036459931310561,1.058853156592918,1.033051371645413,1.062571207143699,1.037931823296638
1.02,1.030908269730445,1.055822624978711,1.030248694515163,1.05966477389862,1.044952854311567,1.03704487911787,1.059079795203124,1.033592528848901,1.062909320420922,1.038517601796928
1.02,1.031136067399556,1.055914068939444,1.030448053293765,1.059794248227486,1.044987233542029,1.037244068952667,1.059156937639645,1.033776800036706,1.063024429240708,1.038717074504137
1.02,1.030116224276146,1.055504618687371,1.029555501246419,1.059214535318953,1.044833108660532,1.036352209475441,1.058811403650558,1.032951710884047,1.062508928475099,1.037823948484156
1.02,1.026781270987358,1.054164723006706,1.02663629376295,1.057317838473081,1.044325433882597,1.033434142449718,1.057678545539637,1.030251632034241,1.060820405741682,1.034901737468549
1.02,1.02466911457333,1.053315450480651,1.024787065705802,1.056115900266754,1.044001113261727,1.031584773037088,1.056958867395153,1.028540047064882,1.059748915125141,1.033049741739228
1.02,1.025790296127977,1.053766321748507,1.025768714220525,1.056753976553141,1.04417353194909,1.032566578297383,1.057341092757154,1.029448739467299,1.060317879752316,1.034032941275607
1.02,1.030171882573573,1.055526968406693,1.029604214501469,1.05924617730552,1.044841533972276,1.036400889155805,1.058830272488253,1.032996747876022,1.062537072806612,1.03787269729525
1.02,1.035192969359673,1.057541227992984,1.033997832248593,1.062098654920482,1.045594823040813,1.040789431659957,1.06052696258283,1.037056036024952,1.065070765653658,1.04226747203313
"@

$rows = $data -split "`r?`n" | Where-Object { $_.Trim() -ne "" }

$rowIndex = 2
foreach ($line in $rows) {
    $fields = $line -split ","
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$rowIndex"
        $ws.Range($cellRef).Value = [double]$fields[$i]
    }
    $rowIndex++
}
